$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "80.402.54"
$ws.Range("E2").Value = "  +5.08%  "

Set-TextValue "D3" "3.211.43"
$ws.Range("E3").Value = "  +4.15%  "

Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.18%  "

Set-TextValue "D5" "210.54"
$ws.Range("E5").Value = "  +6.07%  "

Set-TextValue "D6" "633.48"
$ws.Range("E6").Value = "  +2.74%  "

Set-TextValue "D7" "0.276"
$ws.Range("E7").Value = "  +32.48%  "

Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.04%  "

Set-TextValue "D9" "0.608"
$ws.Range("E9").Value = "  +10.24%  "

Set-TextValue "D10" "3.215.73"
$ws.Range("E10").Value = "  +4.34%  "

Set-TextValue "D11" "0.614"
$ws.Range("E11").Value = "  +39.89%  "

Set-TextValue "D12" "0.0000263"
$ws.Range("E12").Value = "  +35.98%  "

Set-TextValue "D13" "0.166"
$ws.Range("E13").Value = "  +3.39%  "

Set-TextValue "D14" "5.41"
$ws.Range("E14").Value = "  +3.50%  "

Set-TextValue "D15" "3.818.64"
$ws.Range("E15").Value = "  +5.79%  "

Set-TextValue "D16" "32.55"
$ws.Range("E16").Value = "  +12.29%  "

Set-TextValue "D17" "80.610.34"
$ws.Range("E17").Value = "  +5.49%  "

Set-TextValue "D18" "3.224.92"
$ws.Range("E18").Value = "  +5.27%  "

Set-TextValue "D19" "14.59"
$ws.Range("E19").Value = "  +7.44%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "449.26"
$ws.Range("E20").Value = "  +17.81%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "9.37"
$ws.Range("E21").Value = "  +5.06%  "

Set-TextValue "D22" "3.00"
$ws.Range("E22").Value = "  +21.63%  "

Set-TextValue "D23" "5.34"
$ws.Range("E23").Value = "  +21.51%  "

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D24" "6.82"
$ws.Range("E24").Value = "  +5.25%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D25" "3.399.19"
$ws.Range("E25").Value = "  +5.42%  "

Set-TextValue "D26" "77.86"
$ws.Range("E26").Value = "  +7.41%  "

$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D27" "4.79"
$ws.Range("E27").Value = "  +10.49%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D28" "10.98"
$ws.Range("E28").Value = "  +11.23%  "

Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D30" "0.0000125"
$ws.Range("E30").Value = "  +15.77%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "9.29"
$ws.Range("E31").Value = "  +11.93%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D32" "1.00"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D33" "561.41"
$ws.Range("E33").Value = "  +12.44%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D34" "1.51"
$ws.Range("E34").Value = "  +7.95%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D35" "0.155"
$ws.Range("E35").Value = "  +24.10%  "

$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D36" "2.04"
$ws.Range("E36").Value = "  +6.60%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D37" "23.87"
$ws.Range("E37").Value = "  +15.41%  "

$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D38" "0.125"
$ws.Range("E38").Value = "  +21.25%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D39" "0.420"
$ws.Range("E39").Value = "  +10.87%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D40" "0.999"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D41" "166.10"
$ws.Range("E41").Value = "  +2.12%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D42" "5.79"
$ws.Range("E42").Value = "  +13.09%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D43" "20.37"
$ws.Range("E43").Value = "  +1.58%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "192.30"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D45" "1.00"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D46" "1.84"
$ws.Range("E46").Value = "  +12.23%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D47" "2.73"
$ws.Range("E47").Value = "  +12.61%  "

Set-TextValue "D48" "0.803"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D49" "1.34"
$ws.Range("E49").Value = "  +7.82%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D50" "43.76"
$ws.Range("E50").Value = "  +6.16%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D51" "4.34"
$ws.Range("E51").Value = "  +11.99%  "
